# Generate Report for Handoff
# Updates the localization-status workbook to reflect that the
# ca9ac1fd-ed02-42d9-ad23-dfa79d5c9ffc file is now "Ready for handoff"
# (instead of "Handed back: in sync with en-US") and records the new
# handoff timestamp plus an error detail message about a stale handback.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/da7a5ee392da4a0475eea50b307e827772741d48/e2e/ca9ac1fd-ed02-42d9-ad23-dfa79d5c9ffc.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/570bb20fc59a728ab012bdb35579aa2134eeee06/e2e/ca9ac1fd-ed02-42d9-ad23-dfa79d5c9ffc.md."

# ----- Overview sheet: row 3 is the ca9ac1fd-ed02-42d9-ad23-dfa79d5c9ffc.md file -----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-15 09:02:59"

# ----- zh-cn sheet: row 3 is the ca9ac1fd-ed02-42d9-ad23-dfa79d5c9ffc file -----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-08-15 09:02:54"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.1666666666667

# ----- de-de sheet: row 3 is the ca9ac1fd-ed02-42d9-ad23-dfa79d5c9ffc file -----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-08-15 09:02:59"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.1666666666667
